# "integration finale avec BD <3"
# - row 4 : user changed from "monji" to "ahmed"
# - row 5 : user changed from "ahmed" to "karoui"
# - rows 6-12 : reservation timestamps refreshed (new DB export run)
# - rows 13-18 : six brand-new "karoui" / "Confirmé" reservations appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- name corrections on the existing rows ---------------------------------
$ws.Range("A4").Value = "ahmed"
$ws.Range("A5").Value = "karoui"

# --- refreshed timestamps for the already-confirmed rows -------------------
$ws.Range("B6").Value  = "2025-03-06T12:34:49.917466400"
$ws.Range("B7").Value  = "2025-03-06T12:40:07.667345"
$ws.Range("B8").Value  = "2025-03-06T12:44:21.003296400"
$ws.Range("B9").Value  = "2025-03-06T12:44:34.228014700"
$ws.Range("B10").Value = "2025-03-06T12:53:51.661579400"
$ws.Range("B11").Value = "2025-03-06T12:58:30.918849200"
$ws.Range("B12").Value = "2025-03-06T13:01:30.447345700"

# --- append six new reservation rows, cloning row 12's formatting ----------
$newRows = @(
    @{ Row = 13; Date = "2025-03-06T13:02:13.757866900" },
    @{ Row = 14; Date = "2025-03-06T14:12:36.012794600" },
    @{ Row = 15; Date = "2025-03-06 14:46:52" },
    @{ Row = 16; Date = "2025-03-07 08:41:23" },
    @{ Row = 17; Date = "2025-03-07 09:05:09" },
    @{ Row = 18; Date = "2025-03-07 09:26:47" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A12:C12").Copy()
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = "karoui"
    $ws.Range("B" + $r).Value = $item.Date
    $ws.Range("C" + $r).Value = "Confirmé"
}

$excel.CutCopyMode = 0
